$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -10.37309999999999
$ws.Range("D3").Value = -6.994399999999992
$ws.Range("A12").Value = -21.55309999999999
$ws.Range("C14").Value = -12.9669
$ws.Range("D20").Value = -7.558000000000002
$ws.Range("D25").Value = -7.7132
$ws.Range("C26").Value = -12.9296
$ws.Range("A27").Value = -21.755
$ws.Range("D30").Value = -7.272900000000005
$ws.Range("C31").Value = -12.93
$ws.Range("A32").Value = -21.32449999999999
$ws.Range("C35").Value = -11.89930000000001
$ws.Range("A36").Value = -20.0454
$ws.Range("C37").Value = -13.0715
$ws.Range("A38").Value = -19.7109
$ws.Range("D44").Value = -7.449400000000004
$ws.Range("C45").Value = -14.04199999999999
$ws.Range("A46").Value = -21.5122
$ws.Range("D47").Value = -7.4057
$ws.Range("C52").Value = -11.0134
$ws.Range("A54").Value = -21.55799999999998
$ws.Range("A55").Value = -22.48740000000001
$ws.Range("A56").Value = -22.1991
$ws.Range("C57").Value = -14.54399999999998
$ws.Range("D58").Value = -8.230099999999997
$ws.Range("A67").Value = -21.51289999999998
$ws.Range("A69").Value = -21.62899999999997
$ws.Range("A72").Value = -21.46909999999998
$ws.Range("D78").Value = -7.643400000000003
$ws.Range("C81").Value = -13.2328
$ws.Range("A83").Value = -21.41259999999999
$ws.Range("C83").Value = -11.36490000000001
$ws.Range("D84").Value = -8.631000000000004
$ws.Range("A86").Value = -22.37200000000002
$ws.Range("D89").Value = -6.035799999999997
$ws.Range("A91").Value = -21.52470000000001
$ws.Range("D91").Value = -6.043799999999997
$ws.Range("D92").Value = -6.0411
$ws.Range("A93").Value = -21.1824
$ws.Range("D96").Value = -7.435000000000003
$ws.Range("A99").Value = -20.33299999999999
$ws.Range("C100").Value = -12.4915
$ws.Range("C102").Value = -14.53169999999999
$ws.Range("D102").Value = -7.794
